{"js": "// Replace each three-digit-by-one-digit multiplication expression with its\n// updated value, in document order (matches the commit's XML diff).\nconst replacements = [\n  [\"805\u00d78=\", \"153\u00d74=\"],\n  [\"685\u00d76=\", \"154\u00d73=\"],\n  [\"656\u00d73=\", \"946\u00d78=\"],\n  [\"784\u00d73=\", \"417\u00d79=\"],\n  [\"296\u00d78=\", \"858\u00d76=\"],\n  [\"434\u00d78=\", \"567\u00d76=\"],\n  [\"630\u00d79=\", \"410\u00d76=\"],\n  [\"195\u00d79=\", \"733\u00d78=\"],\n  [\"252\u00d77=\", \"504\u00d72=\"],\n  [\"706\u00d77=\", \"476\u00d73=\"],\n  [\"153\u00d75=\", \"401\u00d77=\"],\n  [\"701\u00d77=\", \"349\u00d75=\"],\n  [\"832\u00d78=\", \"862\u00d74=\"],\n  [\"763\u00d77=\", \"209\u00d77=\"],\n  [\"632\u00d73=\", \"878\u00d73=\"],\n  [\"770\u00d78=\", \"723\u00d75=\"],\n  [\"412\u00d72=\", \"157\u00d79=\"],\n  [\"657\u00d77=\", \"962\u00d77=\"],\n  [\"132\u00d75=\", \"556\u00d74=\"],\n  [\"945\u00d77=\", \"617\u00d74=\"],\n  [\"770\u00d79=\", \"492\u00d75=\"],\n  [\"943\u00d75=\", \"174\u00d76=\"],\n  [\"937\u00d77=\", \"143\u00d75=\"],\n  [\"513\u00d74=\", \"607\u00d78=\"],\n  [\"980\u00d78=\", \"145\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  // Replace only the first (and expected only) match for this exact string.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"805\u00d78=\", \"153\u00d74=\"),\n    @(\"685\u00d76=\", \"154\u00d73=\"),\n    @(\"656\u00d73=\", \"946\u00d78=\"),\n    @(\"784\u00d73=\", \"417\u00d79=\"),\n    @(\"296\u00d78=\", \"858\u00d76=\"),\n    @(\"434\u00d78=\", \"567\u00d76=\"),\n    @(\"630\u00d79=\", \"410\u00d76=\"),\n    @(\"195\u00d79=\", \"733\u00d78=\"),\n    @(\"252\u00d77=\", \"504\u00d72=\"),\n    @(\"706\u00d77=\", \"476\u00d73=\"),\n    @(\"153\u00d75=\", \"401\u00d77=\"),\n    @(\"701\u00d77=\", \"349\u00d75=\"),\n    @(\"832\u00d78=\", \"862\u00d74=\"),\n    @(\"763\u00d77=\", \"209\u00d77=\"),\n    @(\"632\u00d73=\", \"878\u00d73=\"),\n    @(\"770\u00d78=\", \"723\u00d75=\"),\n    @(\"412\u00d72=\", \"157\u00d79=\"),\n    @(\"657\u00d77=\", \"962\u00d77=\"),\n    @(\"132\u00d75=\", \"556\u00d74=\"),\n    @(\"945\u00d77=\", \"617\u00d74=\"),\n    @(\"770\u00d79=\", \"492\u00d75=\"),\n    @(\"943\u00d75=\", \"174\u00d76=\"),\n    @(\"937\u00d77=\", \"143\u00d75=\"),\n    @(\"513\u00d74=\", \"607\u00d78=\"),\n    @(\"980\u00d78=\", \"145\u00d74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
